$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (and their relationships) before rebuilding them.
$ws.Hyperlinks.Delete()

# New "Description" column header.
$ws.Range("E1").Value = "Description"

# Update the Image (D) column text for each row to the new picture links.
$abby     = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/Abby%20Doe.png"
$jane     = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/Jane%20Smith.jpg"
$john     = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/John%20Doe.jpg"
$samantha = "https://raw.githubusercontent.com/mostafaalyCS/sheet/refs/heads/main/Samantha%20Black.jpg"
$lorem    = "Lorem Ipsum is simply dummy text of the printing and typesetting industry."

$ws.Range("D2").Value = $abby
$ws.Range("D3").Value = $jane
$ws.Range("D4").Value = $john
$ws.Range("D5").Value = $samantha
$ws.Range("D6").Value = $abby
$ws.Range("D7").Value = $jane
$ws.Range("D8").Value = $john
$ws.Range("D9").Value = $samantha
$ws.Range("D10").Value = $abby
$ws.Range("D11").Value = $jane
$ws.Range("D12").Value = $john
$ws.Range("D13").Value = $samantha

# The "Abby Doe" cells keep their distinct Calibri-based look (no live hyperlink).
$ws.Range("D6").Font.Name = "Calibri"
$ws.Range("D10").Font.Name = "Calibri"

# Fill the new Description column for every data row.
$ws.Range("E2").Value = $lorem
$ws.Range("E3").Value = $lorem
$ws.Range("E4").Value = $lorem
$ws.Range("E5").Value = $lorem
$ws.Range("E6").Value = $lorem
$ws.Range("E7").Value = $lorem
$ws.Range("E8").Value = $lorem
$ws.Range("E9").Value = $lorem
$ws.Range("E10").Value = $lorem
$ws.Range("E11").Value = $lorem
$ws.Range("E12").Value = $lorem
$ws.Range("E13").Value = $lorem

# Re-create the hyperlinks (in the same order they appear in the target file).
# Hyperlinks.Add() resets the cell's font to the default Hyperlink look, so the
# original 20pt sizing is restored immediately afterwards for each cell.
$ws.Hyperlinks.Add($ws.Range("D5"), $samantha)
$ws.Range("D5").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D3"), $jane)
$ws.Range("D3").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D4"), $john)
$ws.Range("D4").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D9"), $samantha)
$ws.Range("D9").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D13"), $samantha)
$ws.Range("D13").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D7"), $jane)
$ws.Range("D7").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D11"), $jane)
$ws.Range("D11").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D8"), $john)
$ws.Range("D8").Font.Size = 20

$ws.Hyperlinks.Add($ws.Range("D12"), $john)
$ws.Range("D12").Font.Size = 20

# Update the active selection.
$ws.Range("G6").Select()
